# "Error Calculations and Plots"
# Remove two data rows (RM 232 and SC 92) and correct a number of
# individual "missing data" cells so the sheet matches the updated
# source table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows for "RM 232" (orig row 26) and "SC 92" (the row
# that slides into position 27 once RM 232 has been removed) ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Apply the scattered cell corrections on the remaining table ---

# Row 5 (RM 14): E column now blank
$ws.Range("E5").ClearContents()

# Row 7 (RM 32): F column now blank
$ws.Range("F7").ClearContents()

# Row 11 (RM 58): E column now has a value
$ws.Range("E11").Value = -7.9

# Row 19 (RM 125): D now has a value, E now blank
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()

# Row 21 (RM 135): D now blank
$ws.Range("D21").ClearContents()

# Row 23 (RM 140): D and E now have values
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7

# Row 24 (RM 142a): F now has a value
$ws.Range("F24").Value = 16.78

# Row 25 (RM 145): E now has a value
$ws.Range("E25").Value = -7.1

# Row 26 (SC 5, after row deletions): C now blank
$ws.Range("C26").ClearContents()

# Row 27 (SC 101, after row deletions): C has a value, D and E now blank
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()

# Row 28 (SC 105, after row deletions): F now has a value
$ws.Range("F28").Value = 17.44

# Row 29 (SC 119, after row deletions): C and E now blank
$ws.Range("C29").ClearContents()
$ws.Range("E29").ClearContents()

# Row 30 (SC 120, after row deletions): F now blank
$ws.Range("F30").ClearContents()

# Row 32 (SC 193, after row deletions): F now blank
$ws.Range("F32").ClearContents()

# Row 33 (SC 232, after row deletions): D and E now have values
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
